# Weekly data refresh for "Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Brócoli"
# Insert 4 new rows at the top of the data block (row 762) which pushes the
# existing rows 762-837 down to 766-841, then populate the 4 new rows with
# the latest week's price records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 762 (shifts old rows 762-837 -> 766-841)
$ws.Rows.Item(762).Resize(4).Insert()

# Values shared by all four new records
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 100112023
$categoria = "Brócoli"
$variedad = "Sin especificar"
$unidad = "`$/unidad"
$kgOUnidades = 1
$clasificacion = "Hortaliza"

# New row 762: Primera, Región Metropolitana
$r = 762
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44769
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 9100
$ws.Cells.Item($r, 11).Value = 600
$ws.Cells.Item($r, 12).Value = 700
$ws.Cells.Item($r, 13).Value = 653
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 653
$ws.Cells.Item($r, 17).Value = $kgOUnidades
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 763: Primera, Región de O'Higgins
$r = 763
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44769
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Primera"
$ws.Cells.Item($r, 10).Value = 3400
$ws.Cells.Item($r, 11).Value = 600
$ws.Cells.Item($r, 12).Value = 600
$ws.Cells.Item($r, 13).Value = 600
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 16).Value = 600
$ws.Cells.Item($r, 17).Value = $kgOUnidades
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 764: Segunda, Región Metropolitana
$r = 764
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44769
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 2700
$ws.Cells.Item($r, 11).Value = 500
$ws.Cells.Item($r, 12).Value = 500
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 500
$ws.Cells.Item($r, 17).Value = $kgOUnidades
$ws.Cells.Item($r, 18).Value = $clasificacion

# New row 765: Segunda, Región de O'Higgins
$r = 765
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = 44769
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $categoriaId
$ws.Cells.Item($r, 7).Value = $categoria
$ws.Cells.Item($r, 8).Value = $variedad
$ws.Cells.Item($r, 9).Value = "Segunda"
$ws.Cells.Item($r, 10).Value = 2800
$ws.Cells.Item($r, 11).Value = 400
$ws.Cells.Item($r, 12).Value = 400
$ws.Cells.Item($r, 13).Value = 400
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 16).Value = 400
$ws.Cells.Item($r, 17).Value = $kgOUnidades
$ws.Cells.Item($r, 18).Value = $clasificacion
